# Update cryptos list (Mon Aug 26 18:50:54 UTC 2024 with GitHub Actions).
# Refreshes per-coin Price / Volume(1h) figures pulled from coinranking.com,
# and picks up three coin-rank swaps that happened between snapshots:
#   Filecoin <-> OKB (rows 41/42), InjectiveProtocol <-> EnergySwap
#   (rows 43/44), VeChain <-> Stellar (rows 50/51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.548.34"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.692.85"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.70"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.88"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("E9").Value = "  -4.04%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  -8.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.171.92"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.51"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.419.96"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.696.28"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.12"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  -5.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.85"
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.35"
$ws.Range("E21").Value = "  -4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.508"
$ws.Range("E23").Value = "  -4.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.00"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.86"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.61"
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "343.99"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.949"
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.10"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.24"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.95"
$ws.Range("E42").Value = "  -5.76%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.45"
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.92"
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0566"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.08"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.27"
$ws.Range("E49").Value = "  -5.30%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0243"
$ws.Range("E51").Value = "  -4.83%  "
